{"js": "// Fill in the header row of the first table with column titles:\n// N\u00b0, Nombre, Carnet, Tel\u00e9fono.\nconst body = context.document.body;\nbody.tables.load(\"items\");\nawait context.sync();\n\nconst table = body.tables.items[0];\nconst headers = [\"N\u00b0\", \"Nombre\", \"Carnet\", \"Tel\u00e9fono\"];\nfor (let col = 0; col < headers.length; col++) {\n  table.getCell(0, col).value = headers[col];\n}\nawait context.sync();\n", "ps1": "# Fill in the header row of the first table with column titles:\n# N\u00b0, Nombre, Carnet, Tel\u00e9fono.\n$d = $word.ActiveDocument\n$table = $d.Tables(1)\n\n$headers = @(\"N\u00b0\", \"Nombre\", \"Carnet\", \"Tel\u00e9fono\")\nfor ($col = 1; $col -le $headers.Length; $col++) {\n    $table.Cell(1, $col).Range.Text = $headers[$col - 1]\n}\n"}
